$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 33: port label becomes "SAN DIEGO AREA TOTALS", species/category label becomes "Totals"
$ws.Range("A33").Value = "SAN DIEGO AREA TOTALS"
$ws.Range("C33").Value = "Totals"

# Update the active selection to match the author's final cursor position
$ws.Range("B30").Select()
